# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.275.65"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.565.16"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.59"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.67"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "3.564.62"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "4.170.18"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.38"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "3.558.67"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "66.334.90"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.53"
$ws.Range("E19").Value = "  +5.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.85"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.55"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.58"
$ws.Range("D25").Value = "3.706.39"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "3.559.82"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.45"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -5.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.84"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.73"
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.01"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").Value = "  +5.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.15"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.44"
$ws.Range("E51").Value = "  +4.40%  "
